$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '72.311.58'
$ws.Cells.Item(2, 5).Value = '  +1.81%  '
$ws.Cells.Item(3, 4).Value = '2.671.70'
$ws.Cells.Item(3, 5).Value = '  +2.12%  '
$ws.Cells.Item(4, 5).Value = '  +0.02%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '602.51'
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  -0.46%  '
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '178.83'
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = '  -0.61%  '
$ws.Cells.Item(7, 5).Value = '  -0.05%  '
$ws.Cells.Item(8, 5).Value = '  -0.57%  '
$ws.Cells.Item(9, 2).Value = 'LidoStakedEther'
$ws.Cells.Item(9, 3).Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Cells.Item(9, 4).Value = '2.671.42'
$ws.Cells.Item(9, 5).Value = '  +2.13%  '
$ws.Cells.Item(10, 2).Value = 'Dogecoin'
$ws.Cells.Item(10, 3).Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '0.174'
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).Value = '  +5.55%  '
$ws.Cells.Item(11, 5).Value = '  +2.15%  '
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '0.358'
$ws.Cells.Item(12, 4).Style = "Normal"
$ws.Cells.Item(12, 5).Value = '  +2.95%  '
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '5.02'
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Cells.Item(13, 5).Value = '  -0.15%  '
$ws.Cells.Item(14, 4).Value = '3.157.61'
$ws.Cells.Item(14, 5).Value = '  +1.90%  '
$ws.Cells.Item(15, 5).Value = '  +3.63%  '
$ws.Cells.Item(16, 4).Value = '72.213.75'
$ws.Cells.Item(16, 5).Value = '  +1.71%  '
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '26.64'
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Cells.Item(18, 4).Value = '2.667.52'
$ws.Cells.Item(18, 5).Value = '  +1.93%  '
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '12.00'
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Cells.Item(19, 5).Value = '  +4.43%  '
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '8.02'
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value = '  +3.20%  '
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '379.42'
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value = '  -0.46%  '
$ws.Cells.Item(22, 5).Value = '  +1.53%  '
$ws.Cells.Item(23, 5).Value = '  +11.55%  '
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '72.53'
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = '  +0.55%  '
$ws.Cells.Item(25, 5).Value = '  -0.82%  '
$ws.Cells.Item(26, 5).Value = '  -0.03%  '
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '10.02'
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(27, 5).Value = '  +4.26%  '
$ws.Cells.Item(28, 4).Value = '2.807.65'
$ws.Cells.Item(28, 5).Value = '  +2.70%  '
$ws.Cells.Item(29, 5).Value = '  +0.13%  '
$ws.Cells.Item(30, 4).Value = '0.0₃0951'
$ws.Cells.Item(30, 5).Value = '  -0.11%  '
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '8.21'
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Cells.Item(31, 5).Value = '  +2.41%  '
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '523.35'
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).Value = '  -1.14%  '
$ws.Cells.Item(33, 5).Value = '  -0.23%  '
$ws.Cells.Item(35, 5).Value = '  -0.08%  '
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '164.82'
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(36, 5).Value = '  -0.33%  '
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '19.60'
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).Value = '  +2.36%  '
$ws.Cells.Item(38, 5).Value = '  +0.86%  '
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '1.40'
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Cells.Item(40, 5).Value = '  -6.33%  '
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '1.85'
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).Value = '  -1.88%  '
$ws.Cells.Item(42, 2).Value = 'RenderToken'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '5.07'
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = '  +1.08%  '
$ws.Cells.Item(43, 2).Value = 'USDe'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '1.00'
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value = '  +0.06%  '
$ws.Cells.Item(44, 5).Value = '  -0.40%  '
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '0.336'
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = '  +1.64%  '
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '39.32'
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value = '  -1.99%  '
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '153.31'
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = '  -0.50%  '
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '3.76'
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(48, 5).Value = '  +3.07%  '
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '0.552'
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = '  +3.93%  '
$ws.Cells.Item(50, 5).Value = '  +3.12%  '
$ws.Cells.Item(51, 5).Value = '  +1.54%  '
